$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new individual stock account rows (NPS stock account additions)
$ws.Range("A11").Value = "hoho222"
$ws.Range("B11").Value = "'009150"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0

$ws.Range("A12").Value = "jiho264"
$ws.Range("B12").Value = "'336370"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0

# Fix selection (circle plot bug fix) - move selection away from C2:C10
$ws.Range("J29").Select()
